$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 7037.4
$ws.Range("I41").Value = 8174.25
$ws.Range("J41").Value = 2490
$ws.Range("K41").Value = 8174.25
$ws.Range("L41").Value = 2490
$ws.Range("M41").Value = -7734.25
$ws.Range("N41").Value = -3370
$ws.Range("H98").Value = 3377.375
$ws.Range("I98").Value = 2503.6667
$ws.Range("K98").Value = 2503.6667
$ws.Range("M98").Value = -1005.6667
$ws.Range("H122").Value = 3377.375
$ws.Range("I122").Value = 2503.6667
$ws.Range("K122").Value = 7511.000100000001
$ws.Range("M122").Value = -5061.000100000001
$ws.Range("H132").Value = 2792.75
$ws.Range("I132").Value = 2118.3242
$ws.Range("K132").Value = 6354.9726
$ws.Range("M132").Value = -3824.9726
$ws.Range("H138").Value = 4838.1875
$ws.Range("I138").Value = 4806.2
$ws.Range("J138").Value = 4841.9067
$ws.Range("K138").Value = 14418.6
$ws.Range("L138").Value = 14525.7201
$ws.Range("M138").Value = -9278.599999999999
$ws.Range("N138").Value = -24805.7201

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11407.388
$ws.Range("I32").Value = 11233.583
$ws.Range("K32").Value = 11233.583
$ws.Range("M32").Value = -10946.583
$ws.Range("H61").Value = 10507371
$ws.Range("I61").Value = 22234246
$ws.Range("J61").Value = 912655.4
$ws.Range("K61").Value = 22234246
$ws.Range("L61").Value = 912655.4
$ws.Range("M61").Value = -22234034
$ws.Range("N61").Value = -913079.4
$ws.Range("H63").Value = 4100.75
$ws.Range("I63").Value = 3300
$ws.Range("J63").Value = 6503
$ws.Range("K63").Value = 3300
$ws.Range("L63").Value = 6503
$ws.Range("M63").Value = -2614
$ws.Range("N63").Value = -7875
$ws.Range("H66").Value = 4100.75
$ws.Range("I66").Value = 3300
$ws.Range("J66").Value = 6503
$ws.Range("K66").Value = 16500
$ws.Range("L66").Value = 32515
$ws.Range("M66").Value = -13068
$ws.Range("N66").Value = -39379
$ws.Range("H74").Value = 2887.4375
$ws.Range("I74").Value = 2887.4375
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 2887.4375
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -2013.4375
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 2887.4375
$ws.Range("I77").Value = 2887.4375
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 14437.1875
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -10069.1875
$ws.Range("N77").ClearContents()
$ws.Range("H110").Value = 4489.7427
$ws.Range("I110").Value = 4274.8213
$ws.Range("J110").Value = 5349.4287
$ws.Range("K110").Value = 4274.8213
$ws.Range("L110").Value = 5349.4287
$ws.Range("M110").Value = -2229.8213
$ws.Range("N110").Value = -9439.4287
$ws.Range("H122").Value = 3586.2354
$ws.Range("I122").Value = 3201.1
$ws.Range("K122").Value = 9603.299999999999
$ws.Range("M122").Value = -7153.299999999999
$ws.Range("H132").Value = 2634644.5
$ws.Range("I132").Value = 3483.111
$ws.Range("J132").Value = 9092950
$ws.Range("K132").Value = 10449.333
$ws.Range("L132").Value = 27278850
$ws.Range("M132").Value = -7919.332999999999
$ws.Range("N132").Value = -27283910
$ws.Range("H136").Value = 10507371
$ws.Range("I136").Value = 22234246
$ws.Range("J136").Value = 912655.4
$ws.Range("K136").Value = 66702738
$ws.Range("L136").Value = 2737966.2
$ws.Range("M136").Value = -66700188
$ws.Range("N136").Value = -2743066.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1459.2069
$ws.Range("I80").Value = 1467.1111
$ws.Range("K80").Value = 1467.1111
$ws.Range("M80").Value = -469.1111000000001
$ws.Range("H83").Value = 1459.2069
$ws.Range("I83").Value = 1467.1111
$ws.Range("K83").Value = 7335.5555
$ws.Range("M83").Value = -2343.5555
$ws.Range("H134").Value = 11114439
$ws.Range("I134").Value = 2516.3333
$ws.Range("K134").Value = 7548.999899999999
$ws.Range("M134").Value = -5013.999899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 18521794
$ws.Range("I31").Value = 20002918
$ws.Range("J31").Value = 7730.25
$ws.Range("K31").Value = 20002918
$ws.Range("L31").Value = 7730.25
$ws.Range("M31").Value = -20002623
$ws.Range("N31").Value = -8320.25
$ws.Range("H34").Value = 18521794
$ws.Range("I34").Value = 20002918
$ws.Range("J34").Value = 7730.25
$ws.Range("K34").Value = 20002918
$ws.Range("L34").Value = 7730.25
$ws.Range("M34").Value = -20002716
$ws.Range("N34").Value = -8134.25
$ws.Range("H52").Value = 61927.25
$ws.Range("I52").Value = 25354.5
$ws.Range("J52").Value = 98500
$ws.Range("K52").Value = 25354.5
$ws.Range("L52").Value = 98500
$ws.Range("M52").Value = -25060.5
$ws.Range("N52").Value = -99088
$ws.Range("H62").Value = 18722.076
$ws.Range("J62").Value = 22820.889
$ws.Range("L62").Value = 22820.889
$ws.Range("N62").Value = -24068.889
$ws.Range("H65").Value = 18722.076
$ws.Range("J65").Value = 22820.889
$ws.Range("L65").Value = 114104.445
$ws.Range("N65").Value = -120344.445
$ws.Range("H99").Value = 17759.924
$ws.Range("I99").Value = 11997.5
$ws.Range("K99").Value = 11997.5
$ws.Range("M99").Value = -10499.5
$ws.Range("H122").Value = 2880.111
$ws.Range("I122").Value = 2802.2
$ws.Range("J122").Value = 2977.5
$ws.Range("K122").Value = 8406.599999999999
$ws.Range("L122").Value = 8932.5
$ws.Range("M122").Value = -5956.599999999999
$ws.Range("N122").Value = -13832.5
$ws.Range("H126").Value = 17759.924
$ws.Range("I126").Value = 11997.5
$ws.Range("K126").Value = 35992.5
$ws.Range("M126").Value = -33522.5
$ws.Range("H132").Value = 1925.0588
$ws.Range("I132").Value = 1781.862
$ws.Range("J132").Value = 2755.6
$ws.Range("K132").Value = 5345.586
$ws.Range("L132").Value = 8266.799999999999
$ws.Range("M132").Value = -2815.586
$ws.Range("N132").Value = -13326.8
$ws.Range("H133").Value = 84133.336
$ws.Range("J133").Value = 84133.336
$ws.Range("L133").Value = 84133.336
$ws.Range("N133").Value = -89193.336
$ws.Range("H134").Value = 1496.1111
$ws.Range("I134").Value = 1496.1111
$ws.Range("K134").Value = 4488.3333
$ws.Range("M134").Value = -1953.3333
$ws.Range("H137").Value = 57499.5
$ws.Range("J137").Value = 57499.5
$ws.Range("L137").Value = 57499.5
$ws.Range("N137").Value = -67699.5
$ws.Range("H141").Value = 282866.34
$ws.Range("J141").Value = 321876.34
$ws.Range("L141").Value = 321876.34
$ws.Range("N141").Value = -332236.34

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 8943.166999999999
$ws.Range("I130").Value = 7830
$ws.Range("J130").Value = 9499.75
$ws.Range("K130").Value = 23490
$ws.Range("L130").Value = 28499.25
$ws.Range("M130").Value = -18470
$ws.Range("N130").Value = -38539.25
$ws.Range("H131").Value = 3535.0715
$ws.Range("J131").Value = 3735.6
$ws.Range("L131").Value = 11206.8
$ws.Range("N131").Value = -21286.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 35646
$ws.Range("I53").Value = 28469.5
$ws.Range("K53").Value = 28469.5
$ws.Range("M53").Value = -27838.5
$ws.Range("H80").Value = 3809.5334
$ws.Range("I80").Value = 2461.4443
$ws.Range("J80").Value = 5831.6665
$ws.Range("K80").Value = 2461.4443
$ws.Range("L80").Value = 5831.6665
$ws.Range("M80").Value = -1463.4443
$ws.Range("N80").Value = -7827.6665
$ws.Range("H83").Value = 3809.5334
$ws.Range("I83").Value = 2461.4443
$ws.Range("J83").Value = 5831.6665
$ws.Range("K83").Value = 12307.2215
$ws.Range("L83").Value = 29158.3325
$ws.Range("M83").Value = -7315.2215
$ws.Range("N83").Value = -39142.3325
$ws.Range("H122").Value = 3703.0527
$ws.Range("I122").Value = 3849.7144
$ws.Range("J122").Value = 3292.4
$ws.Range("K122").Value = 11549.1432
$ws.Range("L122").Value = 9877.200000000001
$ws.Range("M122").Value = -9099.143199999999
$ws.Range("N122").Value = -14777.2
$ws.Range("H132").Value = 3451972.2
$ws.Range("I132").Value = 3718.923
$ws.Range("K132").Value = 11156.769
$ws.Range("M132").Value = -8626.769

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 17877096
$ws.Range("I100").Value = 2047.5
$ws.Range("J100").Value = 31283382
$ws.Range("K100").Value = 2047.5
$ws.Range("L100").Value = 31283382
$ws.Range("M100").Value = -1506.5
$ws.Range("N100").Value = -31284464
$ws.Range("H122").Value = 3484
$ws.Range("I122").Value = 3447.524
$ws.Range("K122").Value = 10342.572
$ws.Range("M122").Value = -7892.572
$ws.Range("H132").Value = 4914.6665
$ws.Range("J132").Value = 6205.3076
$ws.Range("L132").Value = 18615.9228
$ws.Range("N132").Value = -23675.9228

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 279630.8
$ws.Range("I132").Value = 1693.8387
$ws.Range("K132").Value = 5081.5161
$ws.Range("M132").Value = -2551.5161
